# Update 2018-05-25#5 - Database script update on db.error.code.xlsx (Sheet1)
#  - Remove the stray blank-message SignIn row (A=1904).
#  - Replace the stray blank-message CheckAccess row (A=2303) with the new
#    GetAccessUser (4 rows) and SignOut (2 rows) error-code entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Delete the old row 214 (1904 / <blank> / SignIn) - everything below
#    shifts up by one row.
$ws.Rows("214").Delete()

# 2) After the shift above, the row that used to be 2303/CheckAccess
#    (with a blank message) now lives at row 238. Make room for 6 rows
#    (it becomes 6 data rows instead of 1) by inserting 5 blank rows
#    right after it, then fill all six rows in with their new content.
$ws.Rows("239:243").Insert()

$ws.Cells.Item(238, 1).Value = 2303
$ws.Cells.Item(238, 2).Value = "Lang Id cannot be null or empty string."
$ws.Cells.Item(238, 3).Value = "GetAccessUser"

$ws.Cells.Item(239, 1).Value = 2304
$ws.Cells.Item(239, 2).Value = "Lang Id not found."
$ws.Cells.Item(239, 3).Value = "GetAccessUser"

$ws.Cells.Item(240, 1).Value = 2305
$ws.Cells.Item(240, 2).Value = "Access Id cannot be null or empty string."
$ws.Cells.Item(240, 3).Value = "GetAccessUser"

$ws.Cells.Item(241, 1).Value = 2306
$ws.Cells.Item(241, 2).Value = "Access Id not found."
$ws.Cells.Item(241, 3).Value = "GetAccessUser"

$ws.Cells.Item(242, 1).Value = 2307
$ws.Cells.Item(242, 2).Value = "Access Id cannot be null or empty string."
$ws.Cells.Item(242, 3).Value = "SignOut"

$ws.Cells.Item(243, 1).Value = 2308
$ws.Cells.Item(243, 2).Value = "Access Id not found."
$ws.Cells.Item(243, 3).Value = "SignOut"

# 3) Match the author's final selection (B244) recorded in the sheet view.
$ws.Range("B244").Select()
